$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated attendance data for rows 4-9 (columns B=Nama Pegawai, C=Datang, D=Pulang)
$ws.Range("B4").Value = "Yudha Subakti - 340059291"
$ws.Range("C4").Value = "7:02"
$ws.Range("D4").Value = "18:20"

$ws.Range("B5").Value = "Irwan Kurniawan-340016236"
$ws.Range("C5").Value = "7:05"
$ws.Range("D5").Value = "16:52"

$ws.Range("B6").Value = "Agus Santoso"
$ws.Range("C6").Value = "7:06"
$ws.Range("D6").Value = "16:54"

$ws.Range("B7").Value = "MAYA NOVITA SARI - 58827"
$ws.Range("C7").Value = "7:14"
$ws.Range("D7").Value = "18:16"

$ws.Range("B8").Value = "Reny Anggraeni - 34005929"
$ws.Range("C8").Value = "7:17"
$ws.Range("D8").Value = "18:02"

$ws.Range("B9").Value = "Sri Pura - 3400013224"
$ws.Range("C9").Value = "7:33"
$ws.Range("D9").Value = "16:17"

# Row 10 (formerly "7 / Ananto Yanuar / 16:41") is emptied out entirely
$ws.Range("A10").Value = $null
$ws.Range("B10").Value = $null
$ws.Range("C10").Value = $null
$ws.Range("D10").Value = $null

# Remove the now-unused trailing blank row 12 (shifts everything below row 11 up)
$ws.Rows.Item(12).Delete()

# Re-create the merged cells that previously lived on row 12, now on row 11
$ws.Range("E11:F11").Merge()
$ws.Range("G11:H11").Merge()

# Merging recalculates borders on the merged range; restore the original
# (unmerged-row) border/number-format look by copying formats from row 10,
# which still has the untouched original formatting for these columns.
$ws.Range("E10:H10").Copy()
$ws.Range("E11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
